$d = $word.ActiveDocument

# 1. Update the letter date from September 19, 2025 to September 21, 2025
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing address paragraph "969 Story Road, San Jose CA 95122"
#    (the one right after "nan Ho Family Living Trust", not the one that
#    lives inside the table further down) into two separate paragraphs:
#    "969 Story Road" and "San Jose, CA 95122".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "969 Story Road, San Jose CA 95122`r") {
        $savedAscii = $para.Range.Font.NameAscii
        $savedOther = $para.Range.Font.NameOther
        $savedBi    = $para.Range.Font.NameBi
        $savedSize  = $para.Range.Font.Size
        $savedSizeBi = $para.Range.Font.SizeBi

        $para.Range.Find.Execute("969 Story Road, San Jose CA 95122", $true, $false, $false, $false, $false, $true, 1, $false, "969 Story Road^pSan Jose, CA 95122", 2) | Out-Null

        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Font.NameAscii = $savedAscii
        $newPara.Range.Font.NameOther = $savedOther
        $newPara.Range.Font.NameBi    = $savedBi
        $newPara.Range.Font.Size      = $savedSize
        $newPara.Range.Font.SizeBi    = $savedSizeBi

        break
    }
}

# 3. Remove the empty "No Spacing" paragraph that immediately follows the
#    "Board of Directors" signature line.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Board of Directors*") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Style.NameLocal -eq "No Spacing" -and $next.Range.Text -eq "`r") {
            $next.Range.Delete()
        }
        break
    }
}
